$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append (rows 58-60), mirroring the existing table layout
# (columns A-T) used throughout the sheet.
$newRows = @(
    @{ A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Primera";  M=16; N=350000; O=350000; P=350000; Q="`$/bins (450 kilos)"; R="Región Metropolitana";  S=778; T=450 },
    @{ A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Primera";  M=10; N=320000; O=320000; P=320000; Q="`$/bins (450 kilos)"; R="Región de O'Higgins"; S=711; T=450 },
    @{ A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Segunda"; M=8;  N=270000; O=270000; P=270000; Q="`$/bins (450 kilos)"; R="Región de O'Higgins"; S=600; T=450 }
)

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$startRow = 58
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    foreach ($col in $colOrder) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
    # Column D carries the same date number format as the rest of the table.
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "Added rows 58-60"
